$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new inventory item as row 15
$ws.Range("A15").Value = "93J10R"
$ws.Range("B15").Value = "Cilindro o tambor de imagen de fotocopiadora RICOH"
$ws.Range("C15").Value = "1013 1515 MP201 MP301"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 150000
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0
$ws.Range("H15").Formula = "=(E15-D15)*G15"
$ws.Range("I15").Formula = "=D15*F15"
$ws.Range("J15").Value = 0
